# Refresh the Price (D) and Volume(1h) (E) columns on the crypto table
# with the latest scraped figures (GitHub Actions cron update).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 2; Col = 4; Value = '62.491.31'; Text = $true }
    @{ Row = 2; Col = 5; Value = '  -1.44%  '; Text = $false }
    @{ Row = 3; Col = 4; Value = '3.009.43'; Text = $true }
    @{ Row = 3; Col = 5; Value = '  -1.68%  '; Text = $false }
    @{ Row = 4; Col = 5; Value = '  -0.05%  '; Text = $false }
    @{ Row = 5; Col = 4; Value = '584.11'; Text = $true }
    @{ Row = 5; Col = 5; Value = '  -1.46%  '; Text = $false }
    @{ Row = 6; Col = 4; Value = '146.14'; Text = $true }
    @{ Row = 6; Col = 5; Value = '  -5.05%  '; Text = $false }
    @{ Row = 7; Col = 5; Value = '  -0.01%  '; Text = $false }
    @{ Row = 8; Col = 5; Value = '  -2.88%  '; Text = $false }
    @{ Row = 9; Col = 4; Value = '3.005.62'; Text = $true }
    @{ Row = 9; Col = 5; Value = '  -1.72%  '; Text = $false }
    @{ Row = 10; Col = 5; Value = '  -4.06%  '; Text = $false }
    @{ Row = 11; Col = 4; Value = '5.81'; Text = $true }
    @{ Row = 11; Col = 5; Value = '  -0.10%  '; Text = $false }
    @{ Row = 12; Col = 5; Value = '  +1.75%  '; Text = $false }
    @{ Row = 13; Col = 5; Value = '  -3.18%  '; Text = $false }
    @{ Row = 14; Col = 4; Value = '34.66'; Text = $true }
    @{ Row = 14; Col = 5; Value = '  -5.72%  '; Text = $false }
    @{ Row = 15; Col = 5; Value = '  +2.04%  '; Text = $false }
    @{ Row = 16; Col = 4; Value = '3.501.92'; Text = $true }
    @{ Row = 16; Col = 5; Value = '  -1.80%  '; Text = $false }
    @{ Row = 17; Col = 4; Value = '7.08'; Text = $true }
    @{ Row = 17; Col = 5; Value = '  -1.41%  '; Text = $false }
    @{ Row = 18; Col = 4; Value = '62.448.07'; Text = $true }
    @{ Row = 18; Col = 5; Value = '  -1.44%  '; Text = $false }
    @{ Row = 19; Col = 4; Value = '3.008.57'; Text = $true }
    @{ Row = 19; Col = 5; Value = '  -1.83%  '; Text = $false }
    @{ Row = 20; Col = 4; Value = '458.69'; Text = $true }
    @{ Row = 20; Col = 5; Value = '  -5.57%  '; Text = $false }
    @{ Row = 21; Col = 4; Value = '13.96'; Text = $true }
    @{ Row = 21; Col = 5; Value = '  -2.78%  '; Text = $false }
    @{ Row = 22; Col = 5; Value = '  -2.89%  '; Text = $false }
    @{ Row = 23; Col = 5; Value = '  -1.89%  '; Text = $false }
    @{ Row = 24; Col = 4; Value = '81.51'; Text = $true }
    @{ Row = 24; Col = 5; Value = '  -0.93%  '; Text = $false }
    @{ Row = 25; Col = 4; Value = '12.34'; Text = $true }
    @{ Row = 25; Col = 5; Value = '  -4.42%  '; Text = $false }
    @{ Row = 26; Col = 5; Value = '  -9.09%  '; Text = $false }
    @{ Row = 27; Col = 4; Value = '10.01'; Text = $true }
    @{ Row = 27; Col = 5; Value = '  -6.17%  '; Text = $false }
    @{ Row = 28; Col = 4; Value = '1.00'; Text = $true }
    @{ Row = 28; Col = 5; Value = '  +0.20%  '; Text = $false }
    @{ Row = 29; Col = 5; Value = '  -0.08%  '; Text = $false }
    @{ Row = 30; Col = 4; Value = '2.62'; Text = $true }
    @{ Row = 30; Col = 5; Value = '  -2.69%  '; Text = $false }
    @{ Row = 31; Col = 5; Value = '  -4.80%  '; Text = $false }
    @{ Row = 32; Col = 5; Value = '  -5.93%  '; Text = $false }
    @{ Row = 33; Col = 4; Value = '28.10'; Text = $true }
    @{ Row = 33; Col = 5; Value = '  +2.03%  '; Text = $false }
    @{ Row = 34; Col = 5; Value = '  -2.45%  '; Text = $false }
    @{ Row = 35; Col = 4; Value = '0.0₃0810'; Text = $true }
    @{ Row = 35; Col = 5; Value = '  -1.66%  '; Text = $false }
    @{ Row = 36; Col = 4; Value = '1.03'; Text = $true }
    @{ Row = 36; Col = 5; Value = '  -3.18%  '; Text = $false }
    @{ Row = 37; Col = 4; Value = '5.77'; Text = $true }
    @{ Row = 37; Col = 5; Value = '  -3.67%  '; Text = $false }
    @{ Row = 38; Col = 5; Value = '  -5.55%  '; Text = $false }
    @{ Row = 39; Col = 4; Value = '50.31'; Text = $true }
    @{ Row = 39; Col = 5; Value = '  -0.54%  '; Text = $false }
    @{ Row = 40; Col = 5; Value = '  -1.72%  '; Text = $false }
    @{ Row = 41; Col = 4; Value = '2.91'; Text = $true }
    @{ Row = 41; Col = 5; Value = '  -13.28%  '; Text = $false }
    @{ Row = 42; Col = 5; Value = '  +5.28%  '; Text = $false }
    @{ Row = 43; Col = 4; Value = '390.04'; Text = $true }
    @{ Row = 43; Col = 5; Value = '  -11.40%  '; Text = $false }
    @{ Row = 44; Col = 5; Value = '  -1.85%  '; Text = $false }
    @{ Row = 45; Col = 5; Value = '  -7.66%  '; Text = $false }
    @{ Row = 46; Col = 4; Value = '2.729.76'; Text = $true }
    @{ Row = 46; Col = 5; Value = '  -3.97%  '; Text = $false }
    @{ Row = 47; Col = 4; Value = '37.54'; Text = $true }
    @{ Row = 47; Col = 5; Value = '  -3.33%  '; Text = $false }
    @{ Row = 48; Col = 4; Value = '129.35'; Text = $true }
    @{ Row = 48; Col = 5; Value = '  -0.40%  '; Text = $false }
    @{ Row = 50; Col = 5; Value = '  -0.88%  '; Text = $false }
    @{ Row = 51; Col = 5; Value = '  -1.91%  '; Text = $false }
)

foreach ($u in $updates) {
    $value = $u.Value
    if ($u.Text) {
        # Leading apostrophe forces Excel to store/keep this as text
        $value = "'" + $value
    }
    $ws.Cells.Item($u.Row, $u.Col).Value = $value
}
